# Updates cryptos list values (Price / Volume(1h) columns), and shifts rows
# 34-51 to reflect inserted "HuobiToken" row pushing subsequent coins down
# by one position and appending "Decentraland" as the new final row.
#
# Numeric-looking text values (e.g. "0.9995", "1.000") are written with
# a temporary Text number format so Excel keeps them as literal strings
# (matching the original inlineStr cell content) instead of silently
# re-parsing them as numbers. The cell style is reset back to Normal
# immediately afterward so no stray "s" attribute is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.511.14'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.728.48'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9995'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.45%  '
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4809'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.95%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2667'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.08%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').Value = '1.724.56'
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07153'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.64'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6167'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.522'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.18'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = '26.515.77'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9998'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006934'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.66'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.57%  '
$ws.Range('D21').Value = '1.946.58'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.530'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.963'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.286'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.35'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('E27').Value = '  +1.58%  '
$ws.Range('E28').Value = '  -0.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.84'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.978'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08027'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.49%  '
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('E33').Value = '  +2.19%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.616'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6360'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.43%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9929'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.43%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9325'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.11%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.100'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +10.15%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.423'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('B40').Value = 'Quant'
$ws.Range('C40').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '104.99'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -9.02%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.002'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01502'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.51%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.591'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.69%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3905'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.19%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.915'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.38%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1185'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.57%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05332'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.75%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.94'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.825'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.71%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.271'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.24%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3432'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.19%  '
